# Update cryptocurrency price/volume data per Jan 7 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "261.37"
Set-TextValue $ws.Range("E2") "0.87%"

Set-TextValue $ws.Range("D3") "27.10"
Set-TextValue $ws.Range("E3") "0.44%"

Set-TextValue $ws.Range("D4") "4.724"
Set-TextValue $ws.Range("E4") "0.96%"

Set-TextValue $ws.Range("D5") "0.06198"
Set-TextValue $ws.Range("E5") "2.55%"

Set-TextValue $ws.Range("D6") "6.716"
Set-TextValue $ws.Range("E6") "0.62%"

Set-TextValue $ws.Range("D7") "0.8514"
Set-TextValue $ws.Range("E7") "-0.90%"

Set-TextValue $ws.Range("D8") "0.9107"
Set-TextValue $ws.Range("E8") "-0.72%"

Set-TextValue $ws.Range("D9") "0.1407"
Set-TextValue $ws.Range("E9") "0.85%"

Set-TextValue $ws.Range("D10") "0.04806"
Set-TextValue $ws.Range("E10") "-8.34%"

Set-TextValue $ws.Range("D11") "0.07093"
Set-TextValue $ws.Range("E11") "0.11%"

Set-TextValue $ws.Range("D12") "0.03169"
Set-TextValue $ws.Range("E12") "3.45%"

Set-TextValue $ws.Range("D13") "0.09059"
Set-TextValue $ws.Range("E13") "-0.79%"

Set-TextValue $ws.Range("D14") "0.001541"
Set-TextValue $ws.Range("E14") "-0.02%"

Set-TextValue $ws.Range("D15") "0.0006160"
Set-TextValue $ws.Range("E15") "1.54%"

Set-TextValue $ws.Range("D16") "0.006135"
Set-TextValue $ws.Range("E16") "0.83%"

Set-TextValue $ws.Range("D17") "3.466"
Set-TextValue $ws.Range("E17") "-0.08%"

Set-TextValue $ws.Range("D18") "3.172"
Set-TextValue $ws.Range("E18") "0.00%"

Set-TextValue $ws.Range("E19") "-0.35%"

Set-TextValue $ws.Range("E20") "-0.67%"

Set-TextValue $ws.Range("E21") "1.00%"

Set-TextValue $ws.Range("D22") "4.108"
Set-TextValue $ws.Range("E22") "-0.05%"

Set-TextValue $ws.Range("D23") "0.04228"
Set-TextValue $ws.Range("E23") "-0.38%"

Set-TextValue $ws.Range("D24") "0.001218"
Set-TextValue $ws.Range("E24") "0.05%"

Set-TextValue $ws.Range("D25") "0.004123"
Set-TextValue $ws.Range("E25") "2.57%"

Set-TextValue $ws.Range("E26") "0.09%"

Set-TextValue $ws.Range("D40") "0.03916"
Set-TextValue $ws.Range("E40") "1.49%"

Set-TextValue $ws.Range("D41") "0.1113"
Set-TextValue $ws.Range("E41") "-0.18%"

Set-TextValue $ws.Range("D42") "0.004134"
Set-TextValue $ws.Range("E42") "2.71%"

Set-TextValue $ws.Range("E43") "-0.73%"

Set-TextValue $ws.Range("D44") "0.01347"
Set-TextValue $ws.Range("E44") "-11.51%"

Set-TextValue $ws.Range("E45") "-0.26%"

Set-TextValue $ws.Range("E46") "0.09%"

Set-TextValue $ws.Range("D47") "0.03591"
Set-TextValue $ws.Range("E47") "-34.17%"

Set-TextValue $ws.Range("E48") "28.70%"

Set-TextValue $ws.Range("E49") "0.09%"

Set-TextValue $ws.Range("E50") "0.09%"
